$d = $word.ActiveDocument

$arrow = [char]0x2192
$old = "Routine " + $arrow + " interrupt.c"
$new = "Routine " + $arrow + " interrupt.c^pWrapepr " + $arrow + " libc.c"

$d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)

# Move the _GoBack bookmark to the end of the document (end of inserted text)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$endRange = $d.Content
$endRange.Collapse(0)
$d.Bookmarks.Add("_GoBack", $endRange)
